# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" timestamps and the "Priority"
# ("ht" = human translation) flags for the rows that were freshly
# handed off (rows 8, 9, 10, 11, 13, 14 on each localized-language sheet),
# and mirrors the corresponding "Latest HO Xliff Generate Date" on the
# Overview sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$rows = @(8, 9, 10, 11, 13, 14)

foreach ($r in $rows) {
    # Overview sheet: "Latest HO Xliff Generate Date" column (G)
    $wsOverview.Range("G$r").Value = "2016-08-18 16:19:50"

    # zh-cn sheet: "Priority" column (E) and "Latest Handoff Datetime" column (H)
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-18 16:19:44"

    # de-de sheet: "Priority" column (E) and "Latest Handoff Datetime" column (H)
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-18 16:19:50"
}
